$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need an explicit text
# format before assignment, otherwise Excel auto-converts them to numbers
# (the source table stores them as text with trailing padding spaces).
$numericLookingRefs = @("B4","C4","D4","E4","F4","G4","C6","D6","E6","F6","G6","C7","D7","E7","F7","G7","B8","C8","D8","E8","F8","G8","B9","C9","D9","E9","F9","G9","B10","C10","D10","E10","F10","G10","B14","C14","E14","F14","G14","B16","C16","D16","F16","G16","D18","D20","E20","F20","G20","D22","E22","D28","E28","F28","G28")
foreach ($ref in $numericLookingRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated model values (regression output refresh).
$ws.Range("B4").Value = "1808            "
$ws.Range("C4").Value = "1808            "
$ws.Range("D4").Value = "1808            "
$ws.Range("E4").Value = "1808            "
$ws.Range("F4").Value = "1808            "
$ws.Range("G4").Value = "1808            "
$ws.Range("C6").Value = "0.2404          "
$ws.Range("D6").Value = "0.6574          "
$ws.Range("E6").Value = "0.6670          "
$ws.Range("F6").Value = "0.6628          "
$ws.Range("G6").Value = "0.6650          "
$ws.Range("C7").Value = "0.2404          "
$ws.Range("D7").Value = "0.6574          "
$ws.Range("E7").Value = "0.6670          "
$ws.Range("F7").Value = "0.6628          "
$ws.Range("G7").Value = "0.6650          "
$ws.Range("B8").Value = "-0.1076         "
$ws.Range("C8").Value = "-0.1169         "
$ws.Range("D8").Value = "0.9153          "
$ws.Range("E8").Value = "0.8731          "
$ws.Range("F8").Value = "0.8928          "
$ws.Range("G8").Value = "0.8828          "
$ws.Range("B9").Value = "0.0058          "
$ws.Range("C9").Value = "0.0014          "
$ws.Range("D9").Value = "0.8654          "
$ws.Range("E9").Value = "0.8327          "
$ws.Range("F9").Value = "0.8485          "
$ws.Range("G9").Value = "0.8406          "
$ws.Range("B10").Value = "135.56          "
$ws.Range("C10").Value = "94.011          "
$ws.Range("D10").Value = "426.88          "
$ws.Range("E10").Value = "395.95          "
$ws.Range("F10").Value = "388.56          "
$ws.Range("G10").Value = "392.39          "
$ws.Range("B12").Value = "0.0243***       "
$ws.Range("C12").Value = "0.0245***       "
$ws.Range("D12").Value = "0.0224***       "
$ws.Range("E12").Value = "0.0219***       "
$ws.Range("F12").Value = "0.0221***       "
$ws.Range("G12").Value = "0.0220***       "
$ws.Range("B13").Value = "(0.0029)        "
$ws.Range("C13").Value = "(0.0029)        "
$ws.Range("D13").Value = "(0.0009)        "
$ws.Range("E13").Value = "(0.0011)        "
$ws.Range("F13").Value = "(0.0010)        "
$ws.Range("G13").Value = "(0.0010)        "
$ws.Range("B14").Value = "0.0848          "
$ws.Range("C14").Value = "0.0830          "
$ws.Range("D14").Value = "0.0338*         "
$ws.Range("E14").Value = "0.0091          "
$ws.Range("F14").Value = "0.0177          "
$ws.Range("G14").Value = "0.0127          "
$ws.Range("B15").Value = "(0.0563)        "
$ws.Range("C15").Value = "(0.0543)        "
$ws.Range("E15").Value = "(0.0214)        "
$ws.Range("F15").Value = "(0.0204)        "
$ws.Range("B16").Value = "-0.0557         "
$ws.Range("C16").Value = "-0.0563         "
$ws.Range("D16").Value = "-0.0391         "
$ws.Range("E16").Value = "-0.0507*        "
$ws.Range("F16").Value = "-0.0348         "
$ws.Range("G16").Value = "-0.0408         "
$ws.Range("D17").Value = "(0.0247)        "
$ws.Range("E17").Value = "(0.0260)        "
$ws.Range("F17").Value = "(0.0243)        "
$ws.Range("G17").Value = "(0.0249)        "
$ws.Range("B18").Value = "-0.1649*        "
$ws.Range("C18").Value = "-0.1585*        "
$ws.Range("D18").Value = "-0.0488         "
$ws.Range("E18").Value = "-0.0725**       "
$ws.Range("F18").Value = "-0.0620*        "
$ws.Range("G18").Value = "-0.0674*        "
$ws.Range("B19").Value = "(0.0932)        "
$ws.Range("C19").Value = "(0.0905)        "
$ws.Range("D19").Value = "(0.0318)        "
$ws.Range("E19").Value = "(0.0367)        "
$ws.Range("F19").Value = "(0.0346)        "
$ws.Range("G19").Value = "(0.0359)        "
$ws.Range("B20").Value = "0.3380**        "
$ws.Range("C20").Value = "0.3429***       "
$ws.Range("D20").Value = "0.0855          "
$ws.Range("E20").Value = "0.0561          "
$ws.Range("F20").Value = "0.0667          "
$ws.Range("G20").Value = "0.0607          "
$ws.Range("B21").Value = "(0.1366)        "
$ws.Range("C21").Value = "(0.1301)        "
$ws.Range("D21").Value = "(0.0628)        "
$ws.Range("E21").Value = "(0.0708)        "
$ws.Range("F21").Value = "(0.0648)        "
$ws.Range("G21").Value = "(0.0673)        "
$ws.Range("C22").Value = "-0.0340**       "
$ws.Range("D22").Value = "-0.0093         "
$ws.Range("E22").Value = "-0.0064         "
$ws.Range("D23").Value = "(0.0091)        "
$ws.Range("F23").Value = "(0.0088)        "
$ws.Range("C24").Value = "0.0333**        "
$ws.Range("C25").Value = "(0.0145)        "
$ws.Range("D25").Value = "(0.0080)        "
$ws.Range("D26").Value = "0.7221***       "
$ws.Range("E26").Value = "0.6819***       "
$ws.Range("F26").Value = "0.6998***       "
$ws.Range("G26").Value = "0.6906***       "
$ws.Range("D27").Value = "(0.0612)        "
$ws.Range("E27").Value = "(0.0679)        "
$ws.Range("F27").Value = "(0.0676)        "
$ws.Range("G27").Value = "(0.0684)        "
$ws.Range("D28").Value = "0.0225          "
$ws.Range("E28").Value = "-0.0060         "
$ws.Range("F28").Value = "0.0070          "
$ws.Range("G28").Value = "0.0003          "
$ws.Range("D29").Value = "(0.0196)        "
$ws.Range("E29").Value = "(0.0146)        "
$ws.Range("F29").Value = "(0.0166)        "
$ws.Range("G29").Value = "(0.0156)        "
$ws.Range("E30").Value = "0.0985***       "
$ws.Range("E31").Value = "(0.0339)        "
$ws.Range("F32").Value = "0.0572**        "
$ws.Range("F33").Value = "(0.0260)        "
$ws.Range("G34").Value = "0.0791**        "
$ws.Range("G35").Value = "(0.0310)        "

# Restore default (Normal) styling on every touched cell so only the
# cell *content* changes, matching the source formatting.
foreach ($ref in $numericLookingRefs) {
    $ws.Range($ref).Style = "Normal"
}
